$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.210.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.00%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.362.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.06%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.33%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'542.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.68%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'136.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.20%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.26%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +5.56%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +1.92%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'5.58"
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.68%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.55%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +2.40%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.779.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.00%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'58.164.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.91%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +1.54%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.364.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.57%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +3.22%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'333.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.45%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +2.73%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.31%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.08%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'62.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.63%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.168"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.36%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'8.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.47%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.26%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +3.09%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'172.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.56%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +2.47%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +2.65%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.39%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +12.33%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'18.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D35").Value = "'4.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.96%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.47%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +0.72%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.77%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +0.38%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'146.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.68%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'294.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.16%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +0.94%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +1.99%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.0949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.88%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'19.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.20%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +0.89%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +1.51%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +2.83%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").Value = "'0.386"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.98%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'17.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.66%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'11.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.33%  "
$ws.Range("E51").Style = "Normal"
